$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Step 1: remove the old rows for "* GUI", "* Support for VTK files",
# "* Integration" and "* Meetings" (rows 11-14). Excel will shift the
# Total/@Parsiss/@Home rows up (15->11, 16->12, 17->13) and auto-fix
# the SUM() range and the =C15 reference.
# -----------------------------------------------------------------
$ws.Range("B11:B14").EntireRow.Delete()

# -----------------------------------------------------------------
# Step 2: update the header month label
# -----------------------------------------------------------------
$ws.Range("A5").Value2 = "مرداد 99"

# -----------------------------------------------------------------
# Step 3: update the activity rows (6-10)
# -----------------------------------------------------------------
$ws.Range("B6").Value2 = "* Patients Database"
$ws.Range("C6").Value2 = 9
$ws.Range("E6").Value2 = "• Patients Database"

$ws.Range("B7").Value2 = "* GUI"
$ws.Range("C7").Value2 = 6
$ws.Range("E7").Value2 = "• Online Tracking (Registered)"

$ws.Range("B8").Value2 = "* Registration"
$ws.Range("C8").Value2 = 15
$ws.Range("E8").ClearContents()

$ws.Range("B9").Value2 = "* Online Tracking"
$ws.Range("C9").Value2 = 5
$ws.Range("E9").ClearContents()

$ws.Range("B10").Value2 = "* Meetings & other"
$ws.Range("C10").Value2 = 5

# -----------------------------------------------------------------
# Step 4: fix up the totals block (now rows 11-13)
# -----------------------------------------------------------------
$ws.Range("D12").Value2 = 1
$ws.Range("D13").Formula = "=C11-D12"
$ws.Rows.Item(11).RowHeight = 15.6

# -----------------------------------------------------------------
# Step 5: create the blank spacer row 4 with a medium top border
# -----------------------------------------------------------------
$spacer = $ws.Range("A4:E4")
$spacer.ClearFormats()
$spacer.Borders.Item(8).Weight = -4138

# -----------------------------------------------------------------
# Step 6: lightly touch the trailing blank cells in column E that
# remain from the old layout (E10/E11)
# -----------------------------------------------------------------
$e11 = $ws.Range("E11")
$e11.Font.Bold = $true
$e11.Font.Size = 12
$e11.NumberFormat = "#,##0"

$e10 = $ws.Range("E10")
$e10.Borders.Item(8).Weight = 2

# -----------------------------------------------------------------
# Step 7: update the selection like the original author left it
# -----------------------------------------------------------------
$ws.Range("E11").Select()
